# Update of league bases, rotating data for rows 101-103 (Lithuania A Lyga)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lithuania A Lyga")

# Columns used in each data row (A and C/D/E are unchanged)
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# Snapshot current values for rows 101, 102, 103 before overwriting
$row101 = @{}
$row102 = @{}
$row103 = @{}
foreach ($c in $cols) {
    $row101[$c] = $ws.Range("$c" + "101").Value2
    $row102[$c] = $ws.Range("$c" + "102").Value2
    $row103[$c] = $ws.Range("$c" + "103").Value2
}

# Cyclic rotation: row101 <- old row102, row102 <- old row103, row103 <- old row101
foreach ($c in $cols) {
    $ws.Range("$c" + "101").Value = $row102[$c]
    $ws.Range("$c" + "102").Value = $row103[$c]
    $ws.Range("$c" + "103").Value = $row101[$c]
}
